$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Torsion Check" header in P1 (same style as existing header row)
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "Torsion Check"

# P2 is a blank styled cell (matches the header sub-row formatting)
$ws.Range("O2").Copy()
$ws.Range("P2").PasteSpecial(-4122)

# Row 4
$ws.Range("M4").Value = "2L-T25@100"
$ws.Range("N4").Value = "2L-T25@100"
$ws.Range("O4").Value = "2L-T25@100"
$ws.Range("P4").Value = "OK"

# Row 5
$ws.Range("P5").Value = "Overstressed. Please reassess"

# Row 6
$ws.Range("P6").Value = "Overstressed. Please reassess"

# Row 7
$ws.Range("P7").Value = "OK"

# Row 8
$ws.Range("M8").Value = "4L-T12@100"
$ws.Range("N8").Value = "4L-T12@150"
$ws.Range("O8").Value = "4L-T12@100"
$ws.Range("P8").Value = "OK"

# Row 9
$ws.Range("M9").Value = "4L-T12@200"
$ws.Range("N9").Value = "4L-T12@200"
$ws.Range("O9").Value = "4L-T12@200"
$ws.Range("P9").Value = "OK"

# Row 10
$ws.Range("M10").Value = "4L-T12@150"
$ws.Range("N10").Value = "4L-T12@250"
$ws.Range("O10").Value = "4L-T12@150"
$ws.Range("P10").Value = "OK"

# Row 11
$ws.Range("M11").Value = "4L-T12@150"
$ws.Range("N11").Value = "4L-T12@200"
$ws.Range("O11").Value = "4L-T12@150"
$ws.Range("P11").Value = "OK"

# Row 12
$ws.Range("M12").Value = "4L-T12@150"
$ws.Range("N12").Value = "4L-T12@200"
$ws.Range("O12").Value = "4L-T12@150"
$ws.Range("P12").Value = "OK"

# Row 13
$ws.Range("M13").Value = "4L-T12@200"
$ws.Range("N13").Value = "4L-T12@250"
$ws.Range("O13").Value = "4L-T12@200"
$ws.Range("P13").Value = "OK"

# Row 14
$ws.Range("M14").Value = "4L-T12@150"
$ws.Range("N14").Value = "4L-T12@250"
$ws.Range("O14").Value = "4L-T12@150"
$ws.Range("P14").Value = "OK"

# Row 15
$ws.Range("M15").Value = "4L-T12@100"
$ws.Range("N15").Value = "4L-T12@250"
$ws.Range("O15").Value = "4L-T12@100"
$ws.Range("P15").Value = "OK"

# Row 16
$ws.Range("P16").Value = "Overstressed. Please reassess"

# Row 17
$ws.Range("M17").Value = "4L-T12@200"
$ws.Range("N17").Value = "4L-T12@250"
$ws.Range("O17").Value = "4L-T12@200"
$ws.Range("P17").Value = "OK"

# Row 18
$ws.Range("M18").Value = "4L-T12@100"
$ws.Range("N18").Value = "4L-T12@200"
$ws.Range("O18").Value = "4L-T12@100"
$ws.Range("P18").Value = "OK"

# Row 19
$ws.Range("M19").Value = "4L-T12@100"
$ws.Range("N19").Value = "4L-T12@200"
$ws.Range("O19").Value = "4L-T12@100"
$ws.Range("P19").Value = "OK"

# Row 20
$ws.Range("M20").Value = "4L-T12@150"
$ws.Range("N20").Value = "4L-T12@250"
$ws.Range("O20").Value = "4L-T12@150"
$ws.Range("P20").Value = "OK"

# Row 21
$ws.Range("M21").Value = "4L-T12@250"
$ws.Range("N21").Value = "4L-T12@250"
$ws.Range("O21").Value = "4L-T12@250"
$ws.Range("P21").Value = "OK"

# Row 22
$ws.Range("M22").Value = "4L-T12@250"
$ws.Range("N22").Value = "4L-T12@250"
$ws.Range("O22").Value = "4L-T12@250"
$ws.Range("P22").Value = "OK"

# Row 23
$ws.Range("M23").Value = "4L-T12@250"
$ws.Range("N23").Value = "4L-T12@250"
$ws.Range("O23").Value = "4L-T12@250"
$ws.Range("P23").Value = "OK"

# Row 24
$ws.Range("M24").Value = "4L-T12@250"
$ws.Range("N24").Value = "4L-T12@250"
$ws.Range("O24").Value = "4L-T12@250"
$ws.Range("P24").Value = "OK"

# Row 25
$ws.Range("M25").Value = "4L-T12@250"
$ws.Range("N25").Value = "4L-T12@250"
$ws.Range("O25").Value = "4L-T12@250"
$ws.Range("P25").Value = "OK"

# Row 26
$ws.Range("M26").Value = "4L-T12@250"
$ws.Range("N26").Value = "4L-T12@250"
$ws.Range("O26").Value = "4L-T12@250"
$ws.Range("P26").Value = "OK"

# Row 27
$ws.Range("M27").Value = "4L-T12@250"
$ws.Range("N27").Value = "4L-T12@250"
$ws.Range("O27").Value = "4L-T12@250"
$ws.Range("P27").Value = "OK"

# Row 28
$ws.Range("P28").Value = "OK"

# Row 29
$ws.Range("M29").Value = "4L-T12@100"
$ws.Range("N29").Value = "4L-T12@150"
$ws.Range("O29").Value = "4L-T12@100"
$ws.Range("P29").Value = "OK"

# Row 30
$ws.Range("P30").Value = "Overstressed. Please reassess"

# Row 31
$ws.Range("M31").Value = "4L-T12@250"
$ws.Range("N31").Value = "4L-T12@250"
$ws.Range("O31").Value = "4L-T12@250"
$ws.Range("P31").Value = "OK"

# Row 32
$ws.Range("M32").Value = "4L-T12@250"
$ws.Range("N32").Value = "4L-T12@250"
$ws.Range("O32").Value = "4L-T12@250"
$ws.Range("P32").Value = "OK"

# Row 33
$ws.Range("M33").Value = "4L-T12@100"
$ws.Range("N33").Value = "4L-T12@100"
$ws.Range("O33").Value = "4L-T12@100"
$ws.Range("P33").Value = "OK"

# Row 34
$ws.Range("M34").Value = "4L-T12@150"
$ws.Range("N34").Value = "4L-T12@250"
$ws.Range("O34").Value = "4L-T12@150"
$ws.Range("P34").Value = "OK"

# Row 35
$ws.Range("M35").Value = "4L-T12@100"
$ws.Range("N35").Value = "4L-T12@100"
$ws.Range("O35").Value = "4L-T12@100"
$ws.Range("P35").Value = "OK"

# Row 36
$ws.Range("M36").Value = "4L-T12@150"
$ws.Range("N36").Value = "4L-T12@250"
$ws.Range("O36").Value = "4L-T12@150"
$ws.Range("P36").Value = "OK"

# Row 37
$ws.Range("M37").Value = "4L-T12@250"
$ws.Range("N37").Value = "4L-T12@250"
$ws.Range("O37").Value = "4L-T12@250"
$ws.Range("P37").Value = "OK"

# Row 38
$ws.Range("M38").Value = "4L-T12@250"
$ws.Range("N38").Value = "4L-T12@250"
$ws.Range("O38").Value = "4L-T12@250"
$ws.Range("P38").Value = "OK"

# Row 39
$ws.Range("M39").Value = "4L-T12@150"
$ws.Range("N39").Value = "4L-T12@150"
$ws.Range("O39").Value = "4L-T12@150"
$ws.Range("P39").Value = "OK"

# Row 40
$ws.Range("M40").Value = "4L-T16@150"
$ws.Range("N40").Value = "4L-T12@250"
$ws.Range("O40").Value = "4L-T16@150"
$ws.Range("P40").Value = "OK"

# Row 41
$ws.Range("M41").Value = "4L-T12@100"
$ws.Range("N41").Value = "4L-T12@100"
$ws.Range("O41").Value = "4L-T12@100"
$ws.Range("P41").Value = "OK"

# Row 42
$ws.Range("M42").Value = "4L-T12@150"
$ws.Range("N42").Value = "4L-T12@150"
$ws.Range("O42").Value = "4L-T12@150"
$ws.Range("P42").Value = "OK"

# Row 43
$ws.Range("M43").Value = "4L-T12@100"
$ws.Range("N43").Value = "4L-T12@250"
$ws.Range("O43").Value = "4L-T12@100"
$ws.Range("P43").Value = "OK"

# Row 44
$ws.Range("M44").Value = "4L-T12@150"
$ws.Range("N44").Value = "4L-T12@150"
$ws.Range("O44").Value = "4L-T12@150"
$ws.Range("P44").Value = "OK"

# Row 45
$ws.Range("P45").Value = "Overstressed. Please reassess"

# Row 46
$ws.Range("P46").Value = "Overstressed. Please reassess"

# Row 47
$ws.Range("P47").Value = "Overstressed. Please reassess"

# Row 48
$ws.Range("M48").Value = "4L-T12@100"
$ws.Range("N48").Value = "4L-T12@100"
$ws.Range("O48").Value = "4L-T12@100"
$ws.Range("P48").Value = "OK"

# Row 49
$ws.Range("P49").Value = "Overstressed. Please reassess"

# Row 50
$ws.Range("P50").Value = "Overstressed. Please reassess"

# Row 51
$ws.Range("P51").Value = "Overstressed. Please reassess"

# Row 52
$ws.Range("P52").Value = "Overstressed. Please reassess"

# Row 53
$ws.Range("P53").Value = "Overstressed. Please reassess"

# Row 54
$ws.Range("P54").Value = "Overstressed. Please reassess"

# Row 55
$ws.Range("P55").Value = "Overstressed. Please reassess"

# Row 56
$ws.Range("P56").Value = "Overstressed. Please reassess"

# Row 57
$ws.Range("P57").Value = "Overstressed. Please reassess"

# Row 58
$ws.Range("P58").Value = "Overstressed. Please reassess"

# Row 59
$ws.Range("P59").Value = "Overstressed. Please reassess"
